$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A5").Value = 'Nike SB Zoom Stefan Janoski "Medium Mint"'
$ws.Range("B5").Value = 1
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "$30.00"

$ws.Range("C7").Select()
